$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header cells to use new shared string values.
$ws.Range("A1").Value = "AA"
$ws.Range("B1").Value = "BB"

# Move the active selection to A2 (matches sheetView selection change in diff).
$ws.Range("A2").Select()
